# EST-1134: fleshes out IncomingCharge with fixture and integtest
#
# Reworks the ChargeHierarchy sheet:
#  - header row becomes Name / Parent / AtPath (was Charges / ChargeParent)
#  - a new "FRANCE" top-level row is inserted under the header (row 2)
#  - every charge-type row gets an AtPath of "/FRA" in column C
#  - "LEGAL FEES" is renamed to "LEGAL / BAILIFF FEES"
#  - the active selection moves from B4 to C4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "LEGAL FEES" to "LEGAL / BAILIFF FEES" first (row 8, column A)
# before anything else touches the sheet.
$ws.Range("A8").Value = "LEGAL / BAILIFF FEES"

# New header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Parent"
$ws.Range("C1").Value = "AtPath"

# Full list of charge-type names (column A, rows 3-13), in order.
$chargeTypes = @(
    "PROJECT MANAGEMENT",
    "TAX",
    "WORKS",
    "RELOCATION / DISPOSSESSION INDEMNITY",
    "ARCHITECT / GEOMETRICIAN FEES",
    "LEGAL / BAILIFF FEES",
    "MARKETING",
    "TENANT INSTALLATION WORKS",
    "SECURITY AGENTS",
    "LETTING FEES",
    "OTHER"
)

# Rows 3-13: each charge type, parented under FRANCE, AtPath "/FRA".
# Leading apostrophe => quote-prefixed text entry (AtPath values start
# with "/", which Excel treats as text needing the prefix flag).
for ($i = 0; $i -lt $chargeTypes.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $chargeTypes[$i]
    $ws.Cells.Item($row, 2).Value = "FRANCE"
    $ws.Cells.Item($row, 3).Value = "'/FRA"
}

# Row 2: new top-level "FRANCE" entry, no parent, AtPath "/FRA"
$ws.Range("A2").Value = "FRANCE"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "'/FRA"

# Selection moves to C4 in the saved file
$ws.Range("C4").Select()
